$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.324.01"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").Value = "2.051.88"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.621"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.63%  "
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0769"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.55%  "
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "2.356.25"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.757"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("D17").Value = "2.050.02"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "37.338.62"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("D21").Value = "0.0₃0823"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.70%  "
$ws.Range("E26").Value = "  +4.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.130"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.38%  "
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0621"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.56%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  -2.66%  "
$ws.Range("E40").Value = "  +4.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").Value = "1.486.38"
$ws.Range("E42").Value = "  +2.97%  "
$ws.Range("E43").Value = "  -3.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.41%  "
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.05%  "
$ws.Range("E50").Value = "  -2.49%  "
$ws.Range("D51").Value = "2.242.18"
$ws.Range("E51").Value = "  -1.13%  "
